# "updated reports, fixed table element padding, added relocation button to card"
#
# The "Code refactoring" (row 20 / H20) and "Front-end: Relocation pop-up" (row 18)
# comment text is cleared out as part of the reports update: the two long
# freeform comments tied to the "Code refactoring" task (H20/I20) are removed
# and row 20's "Actual Finish" (G20) is now filled in with the same date as
# its planned Finish (F20), since the task is now complete.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 20 ("Code refactoring"): mark it finished on schedule by copying the
# Finish date/format into Actual Finish, and clear out the now-stale
# freeform comments in H20/I20 (keep their wrap-text style, just blank them).
$ws.Range("F20").Copy()
$ws.Range("G20").PasteSpecial(-4122)   # xlPasteFormats - reuse F20's date style
$ws.Range("G20").Value = $ws.Range("F20").Value2
$ws.Application.CutCopyMode = $false

$ws.Range("H20").ClearContents()
$ws.Range("I20").ClearContents()

# Row 20 no longer needs its tall custom height now that the comments are gone.
$ws.Rows.Item(20).AutoFit()

# Scroll the view down toward the bottom of the table and move the active
# selection from E21 to F21.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F21").Select()

$wb.Save()
